$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.16
$ws.Range("C2").Value = 0.21
$ws.Range("D2").Value = 0.39

$ws.Range("B3").Value = 2.35
$ws.Range("C3").Value = 0.17
$ws.Range("D3").Value = 3

$ws.Range("B4").Value = 0.77
$ws.Range("C4").Value = 0.15
$ws.Range("D4").Value = 0.41
